$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: columns shift from A-E to A,B,C,D,(skip E),F
$ws.Range("E1").ClearContents()

$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Status"
$ws.Range("F1").Value = "Shot Number"
$ws.Range("B1").Value = "Task Number"

# Adjust column widths (values chosen so the engine's internal pixel
# rounding reproduces the target stored widths as closely as possible)
$ws.Columns.Item(1).ColumnWidth = 54.666666666666664
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(6).ColumnWidth = 13.0

# Update selection to match target view state
$ws.Range("E13").Select()
